# "Allow momentary contact setting for pads"
#
# 1) Slide 4 ("Parameter Widgets"): add a new bullet paragraph right after
#    the "Knob resolution" bullet, documenting the Pad caption format and
#    the optional "_m" (momentary contact) suffix. The shape uses
#    <a:spAutoFit/>, so its height grows automatically to fit the new line.
#
# 2) Footer "date" placeholders on the slide master + all slide layouts:
#    refresh the cached datetimeFigureOut display text (cosmetic, an
#    artifact of re-saving the deck on a later date).

$p = $ppt.ActivePresentation

# --- 1. New "Pad" bullet on slide 4 -----------------------------------
$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$total = $tr.Paragraphs().Count
for ($i = 1; $i -le $total; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.IndexOf("Knob resolution") -ge 0) {
        $knobPara = $para
        break
    }
}

$dash = [char]8211
$lq = [char]8220
$rq = [char]8221

$newText = "Pad " + $dash + " Caption format is " + $lq + "Label_OnText_OffText" + $rq + " with optional " + $lq + "_m" + $rq + " appended make it momentary contact"

# Split the new paragraph off from the "Knob resolution" one; it inherits
# that paragraph's bullet/indent formatting automatically.
$knobPara.InsertAfter("`r" + $newText)

# Re-split the middle run ("Label_OnText_OffText") out from its neighbours
# by re-asserting its (unchanged) font size, matching the run boundaries
# produced when the original author's text editor flagged that token.
$total2 = $tr.Paragraphs().Count
for ($i = 1; $i -le $total2; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.IndexOf("Label_OnText_OffText") -ge 0) {
        $padPara = $para
        break
    }
}
$full = $padPara.Text
$idx = $full.IndexOf("Label_OnText_OffText")
$subStart = $padPara.Start + $idx
$sub = $tr.Characters($subStart, 20)
$sub.Font.Size = 12

# --- 2. Refresh cached footer date text -------------------------------
$oldDate = "3/28/2023"
$newDate = "8/23/2025"

$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $sm.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
